$wb = $excel.ActiveWorkbook

$oldId = "47181a61-8ac6-4af1-b013-0f00d62e9ca2"
$newId = "f56b1430-3df6-465d-8415-390da41f3cbd"
$newHash = "330b549635a956dc37620130096bf83f7993e9e8"

# Hyperlink font colour used by the workbook's "HyperLink" cell style
# (#6495ED, underlined) - reapplied after hyperlinks are recreated below.
$hlColor = 15570276

# ---------------------------------------------------------------------------
# Overview sheet - just text / date updates, hyperlink target text changes
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("A2").Value = "$newId.md"
$ov.Range("G2").Value = "2016-08-12 21:13:35"

$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/e90a835daa426dec55fcd05ab54005cb0c2f3b8e/e2e/$newId.md", [Type]::Missing, [Type]::Missing, "e2e\$newId.md")
$ov.Range("B2").Font.Underline = 2
$ov.Range("B2").Font.Color = $hlColor

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("A2").Value = "$newId.md"
$zh.Range("G2").Value = "$newId.$newHash.zh-cn.xlf"
$zh.Range("H2").Value = "2016-08-12 21:13:27"
$zh.Range("I2").Value = ""
$zh.Range("J2").Value = ""
$zh.Range("K2").Value = "0001-01-01 00:00:00"
$zh.Range("I2").Style = "Normal"
$zh.Range("J2").Style = "Normal"

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/e90a835daa426dec55fcd05ab54005cb0c2f3b8e/e2e/$newId.md", [Type]::Missing, [Type]::Missing, "$newId.md")
$zh.Range("A2").Font.Underline = 2
$zh.Range("A2").Font.Color = $hlColor

$zh.Columns.Item(9).ColumnWidth = 17.833333333333332
$zh.Columns.Item(10).ColumnWidth = 20.833333333333332

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("A2").Value = "$newId.md"
$de.Range("G2").Value = "$newId.$newHash.de-de.xlf"
$de.Range("H2").Value = "2016-08-12 21:13:35"
$de.Range("I2").Value = ""
$de.Range("J2").Value = ""
$de.Range("K2").Value = "0001-01-01 00:00:00"
$de.Range("I2").Style = "Normal"
$de.Range("J2").Style = "Normal"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/e90a835daa426dec55fcd05ab54005cb0c2f3b8e/e2e/$newId.md", [Type]::Missing, [Type]::Missing, "$newId.md")
$de.Range("A2").Font.Underline = 2
$de.Range("A2").Font.Color = $hlColor

$de.Columns.Item(9).ColumnWidth = 17.833333333333332
$de.Columns.Item(10).ColumnWidth = 20.833333333333332

Write-Output "done"
